# Trade #1 (MarketMaking) closed at 2026-02-17 19:55:35 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1399.85   # Current Capital
$wsSummary.Range("B4").Value = -0.15     # Total P&L $
$wsSummary.Range("B5").Value = -3        # Total P&L %
$wsSummary.Range("B6").Value = 1         # Total Trades
$wsSummary.Range("B8").Value = 1         # Losing Trades

# ---- Strategy Status sheet (row 5 = MarketMaking) ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 99.84999999999999  # Capital
$wsStatus.Range("D5").Value = 1                  # Trades
$wsStatus.Range("E5").Value = -0.15              # P&L $
$wsStatus.Range("F5").Value = -0.15              # P&L %

# ---- All Trades sheet (row 2 = Trade #1) ----
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("G2").Value = 0.34            # Exit Price
$wsAllTrades.Range("H2").Value = "CLOSED"        # Status
$wsAllTrades.Range("I2").Value = -30.6122        # P&L %
$wsAllTrades.Range("J2").Value = -0.15           # P&L $
$wsAllTrades.Range("K2").Value = 99.84999999999999  # Capital After
$wsAllTrades.Range("P2").Value = "early_exit"    # Exit Reason
$wsAllTrades.Range("Q2").Value = 5.05            # Duration (min)

# ---- MarketMaking sheet (row 2 = Trade #1, mirrors All Trades) ----
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G2").Value = 0.34                   # Exit Price
$wsMM.Range("H2").Value = "CLOSED"               # Status
$wsMM.Range("I2").Value = -30.6122               # P&L %
$wsMM.Range("J2").Value = -0.15                  # P&L $
$wsMM.Range("K2").Value = 99.84999999999999         # Capital After
$wsMM.Range("P2").Value = "early_exit"           # Exit Reason
$wsMM.Range("Q2").Value = 5.05                   # Duration (min)
